$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - copy formatting from row 15 (same column styles) then set values
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "First email via app"
$ws.Range("C16").Value = 44992.729375
$ws.Range("D16").Value = "Wish me luck!"
$ws.Range("E16").Value = "marleevaughn@outlook.com"
$ws.Range("F16").Value = "Marlee Vaughn"
$ws.Range("G16").Value = "duanevaughn@hotmail.com"
$ws.Range("H16").Value = "Duane Vaughn"
$ws.Range("I16").Value = $false

# Row 17 - copy formatting from row 15 (same column styles) then set values
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "First draft"
$ws.Range("C17").Value = 44992.73216435185
$ws.Range("D17").Value = "My first draft via app has been edited ONCE! and has id of 15 and has now been sent"
$ws.Range("E17").Value = "marleevaughn@outlook.com"
$ws.Range("F17").Value = "Marlee Vaughn"
$ws.Range("G17").Value = "kalevaughn@gmail.com"
$ws.Range("H17").Value = "Kale Vaughn"
$ws.Range("I17").Value = $false
